$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$ws.Range("B24").Value = "474cf07313af8791624506f75e5f63ca"
$ws.Range("B34").Value = "416571f9d44722b4fddf9f3176079781"
$ws.Range("B136").Value = "145f6cdd9e574970a49058607a4c57c6"
$ws.Range("B159").Value = "dbfc21f7e94c2499a7e91e097f364003"
$ws.Range("B162").Value = "93193683d95b78daac8c776087db4190"
$ws.Range("B169").Value = "d8e2d3b430620fbcc36650018a5d213d"
$ws.Range("B180").Value = "5ef31b906e5e014b2a712c7917b67d23"
$ws.Range("B213").Value = "db623fc7026992cb80742ec8618477a7"
$ws.Range("B227").Value = "366679d9cd102f7c634ebffd2d642faa"
$ws.Range("B232").Value = "2ad3ae0d1889ca9238638c3c5377ba7a"
$ws.Range("B281").Value = "7f6ab24a2600337270ff3e0396ae3efd"
$ws.Range("B302").Value = "0f1ef506e706195dbd93c49065f789b1"
$ws.Range("B339").Value = "4355b8ccd9f3d91560badc347230afcd"
$ws.Range("B460").Value = "ef3bb11c9a11290215fab20c3653025e"
$ws.Range("B461").Value = "4d3123e6427b8be8bc3bb96a1e1c47a7"
$ws.Range("B478").Value = "0e421a028fe726870a018a31b7132a98"
$ws.Range("B500").Value = "90638a5840cb2ea45547ac598d99705e"
$ws.Range("B501").Value = "10add39a694426657601535a2ecb2c04"
$ws.Range("B502").Value = "81629ac93065ab0b8af54b4a0aeeec72"
$ws.Range("B506").Value = "a2e524582dcc998df58cb03cc9fd0f7d"
$ws.Range("B514").Value = "945b250e5829cb718f5588669ddd991e"
$ws.Range("B517").Value = "d58681c86cbed19c395aab18d70338ab"
$ws.Range("B524").Value = "7a9f409bb9d824128a198556a9c68d64"
$ws.Range("B550").Value = "8aab137630c87b0adee966d8555f7e13"
$ws.Range("B616").Value = "078638d89707ef761041c1aa1f6eb798"
$ws.Range("B627").Value = "0225aa8685f6b6a513936ce0d53587e9"
$ws.Range("B665").Value = "1ba24c61578dfbe6dd75691af4a3de32"
$ws.Range("B666").Value = "eab8275c1ef6f5796f0d9ea05abaa178"
$ws.Range("B680").Value = "c56276a6b66cd48ad0785d014eb73047"
$ws.Range("B685").Value = "ec3c9dece34deb785b2e0c3199198bf1"
$ws.Range("B700").Value = "cf0a52c92f73b57c3c83178f85143e6b"
$ws.Range("B703").Value = "638699b23a2f1096340b7b6ea891a02e"
$ws.Range("B704").Value = "6ce535b8d351b4b4aa19d5896b319229"
$ws.Range("B729").Value = "55111d5fb891311a52db4d42d9478720"
$ws.Range("B742").Value = "fec6821cbe9c86068c0b2ce65f3d2782"
$ws.Range("B819").Value = "ddcecae74f700d34aeb688e4eafe9966"
$ws.Range("B830").Value = "878f501c6fcfbb24100b756563e49341"
$ws.Range("B835").Value = "44a1dc031076aedec8ddf2465a2c79d5"
$ws.Range("B854").Value = "aea50cacf37de8405a6e0d39d5a91d54"

$wb.Save()
